$d = $word.ActiveDocument

# --- Change 1: collapse the multiple runs in the "Then percentage..." paragraph
# into a single run with the concatenated text.
$target = $d.Paragraphs.Item(14).Range
$target.Find.Execute(
    "Then percentage of total employment is number of hours per week divided by 42 (external parameter)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Then percentage of total employment is number of hours per week divided by 42 (external parameter)",
    2
)

# --- Change 2: insert a new list paragraph right after the "For all hourly
# workers, recuperation..." paragraph, with the new bullet text split across
# four runs (mirrors the target OOXML run layout).
$srcPara = $d.Paragraphs.Item(15)
$srcPara.Range.Copy()
$insertionPoint = $d.Range($srcPara.Range.End, $srcPara.Range.End)
$insertionPoint.Paste()

$newPara = $d.Paragraphs.Item(16)
$clearRange = $newPara.Range
$clearRange.End = $clearRange.End - 1
$clearRange.Delete()

$fillPoint = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$runXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Set notice for hourly worker </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>&#8211;</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>according to average wage in last 3 months</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$fillPoint.InsertXML($runXml)
